# "change C to FCU" -- adds a statistics / frequency-distribution block
# below the existing coin-measurement table (MAX/MIN/range/class width,
# plus a grouped frequency table with a FREQUENCY() array formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary stats: MAX / MIN / range (S) / class width (S/5) ---------
$ws.Range("A15").Value = "MAX"
$ws.Range("B15").Formula = "=MAX(G2:G13)"
$ws.Range("B15").NumberFormat = "0.000_ "

$ws.Range("A16").Value = "MIN"
$ws.Range("B16").Formula = "=MIN(G2:G13)"
$ws.Range("B16").NumberFormat = "0.000_ "

$ws.Range("A17").Value = "S"
$ws.Range("B17").Formula = "=B15-B16"
$ws.Range("B17").NumberFormat = "0.000_ "

$ws.Range("A18").Value = "S/5"
$ws.Range("B18").Formula = "=B17/5"
$ws.Range("B18").ClearFormats() | Out-Null

# --- Frequency-distribution table header (same centered look as row 1) --
$ws.Range("A20").Value = "No"
$ws.Range("B20").Value = "group"
$ws.Range("D20").Value = "count"
$ws.Range("E20").Value = "count"
$ws.Range("A20:E20").HorizontalAlignment = -4108
$ws.Range("A20:E20").VerticalAlignment = -4108
$ws.Range("B20:C20").Merge() | Out-Null

# --- Frequency-distribution table body (5 class bins) --------------------
$ws.Range("A21").Value = 1
$ws.Range("B21").Formula = "=B16+0*B17"
$ws.Range("C21").Formula = "=B16+0.2*B17-0.0001"

$ws.Range("A22").Value = 2
$ws.Range("B22").Formula = "=B16+0.2*B17"
$ws.Range("C22").Formula = "=B16+0.4*B17-0.0001"

$ws.Range("A23").Value = 3
$ws.Range("B23").Formula = "=B16+0.4*B17"
$ws.Range("C23").Formula = "=B16+0.6*B17-0.0001"

$ws.Range("A24").Value = 4
$ws.Range("B24").Formula = "=B16+0.6*B17"
$ws.Range("C24").Formula = "=B16+0.8*B17-0.0001"

$ws.Range("A25").Value = 5
$ws.Range("B25").Formula = "=B16+0.8*B17"
$ws.Range("C25").Formula = "=B16+1*B17"

# Array formula spilling the FREQUENCY() counts down D21:D25
$ws.Range("D21:D25").FormulaArray = "=FREQUENCY(G2:G13,C21:C25)"

# The B/C formulas above reference B16/B17 (custom "0.000_ " format), and
# this engine - like Excel's own "format painting by precedent" - copies
# that format onto new, still-General cells that reference them. The
# source workbook keeps this whole block in the plain default format, so
# reset it explicitly.
$ws.Range("B21:C25").ClearFormats() | Out-Null
$ws.Range("D21:D25").ClearFormats() | Out-Null

# Recalculate so every new formula (especially the array formula) carries
# a fresh cached value before we read/inspect anything else.
$excel.Calculate()

# --- Cosmetics -------------------------------------------------------
# Column B now holds decimals, so give it a touch more width (mirrors the
# bestFit column Excel applies automatically after typing new data).
$ws.Columns.Item(2).ColumnWidth = 9.5

# Move the live selection the way the author left it.
$ws.Range("O17").Select() | Out-Null
